$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The hotel cost table's "Total" column header becomes "Total Per Room"
$ws.Range("D1").Value = "Total Per Room"

# Column A was widened slightly (manual resize) while reviewing the results
$ws.Columns("A").ColumnWidth = 31.5

# Final cursor position before the workbook was saved/sent
$ws.Range("G9").Select() | Out-Null
